# Monte Carlo Part Two.docx - remove the "click here" shiny-app hyperlink
# sentence from the opening paragraph.
#
# Before:
#   "...and visualize the results. If you just can’t wait, a link to the
#    final Shiny app is available [here](hyperlink)."
# After:
#   "...and visualize the results. "
#
# Approach: unlink the hyperlink field (turns the "here" run into plain
# text, matching how Word itself removes a hyperlink), then delete the
# now-plain-text sentence "If you just can’t wait, ... available here."
# that trails the sentence we want to keep - leaving the single trailing
# space already present after "results.".

$d = $word.ActiveDocument

# 1) Remove the hyperlink around "here" (unlinks the field; Word keeps the
#    display text "here" as plain text, which we then delete below along
#    with the rest of the now-unwanted sentence).
if ($d.Hyperlinks.Count -gt 0) {
    $d.Hyperlinks.Item(1).Delete()
}

# 2) Delete the trailing sentence, including the final period, while
#    keeping the single space that already follows "results.".
$rng = $d.Content
$found = $rng.Find.Execute(
    "If you just can" + [char]0x2019 + "t wait, a link to the final Shiny app is available here.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Delete()
}
